$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 495.26315
$ws.Range("I9").Value = 607.26666
$ws.Range("K9").Value = 607.26666
$ws.Range("M9").Value = -438.26666
$ws.Range("H12").Value = 5447.846
$ws.Range("I12").Value = 5236.1665
$ws.Range("K12").Value = 5236.1665
$ws.Range("M12").Value = -5066.1665
$ws.Range("H80").Value = 743.4375
$ws.Range("I80").Value = 426.33334
$ws.Range("J80").Value = 933.7
$ws.Range("K80").Value = 1279.00002
$ws.Range("L80").Value = 2801.1
$ws.Range("M80").Value = -281.0000199999999
$ws.Range("N80").Value = -4797.1
$ws.Range("H83").Value = 743.4375
$ws.Range("I83").Value = 426.33334
$ws.Range("J83").Value = 933.7
$ws.Range("K83").Value = 3837.00006
$ws.Range("L83").Value = 8403.300000000001
$ws.Range("M83").Value = 1154.99994
$ws.Range("N83").Value = -18387.3
$ws.Range("H86").Value = 3533
$ws.Range("I86").Value = 4551.4443
$ws.Range("J86").Value = 1699.8
$ws.Range("K86").Value = 4551.4443
$ws.Range("L86").Value = 1699.8
$ws.Range("M86").Value = -3428.4443
$ws.Range("N86").Value = -3945.8
$ws.Range("H87").Value = 79999
$ws.Range("J87").Value = 79999
$ws.Range("L87").Value = 79999
$ws.Range("N87").Value = -82495
$ws.Range("H88").Value = 1509.2727
$ws.Range("J88").Value = 696.5
$ws.Range("L88").Value = 696.5
$ws.Range("N88").Value = -1508.5
$ws.Range("H89").Value = 3533
$ws.Range("I89").Value = 4551.4443
$ws.Range("J89").Value = 1699.8
$ws.Range("K89").Value = 22757.2215
$ws.Range("L89").Value = 8499
$ws.Range("M89").Value = -17141.2215
$ws.Range("N89").Value = -19731
$ws.Range("H90").Value = 79999
$ws.Range("J90").Value = 79999
$ws.Range("L90").Value = 239997
$ws.Range("N90").Value = -252477
$ws.Range("H91").Value = 1509.2727
$ws.Range("J91").Value = 696.5
$ws.Range("L91").Value = 696.5
$ws.Range("N91").Value = -3504.5
$ws.Range("H98").Value = 1117
$ws.Range("I98").Value = 579.36365
$ws.Range("K98").Value = 579.36365
$ws.Range("M98").Value = 918.63635
$ws.Range("H99").Value = 344.33334
$ws.Range("J99").Value = 468
$ws.Range("L99").Value = 1404
$ws.Range("N99").Value = -4400
$ws.Range("H100").Value = 2382.0715
$ws.Range("I100").Value = 2590.818
$ws.Range("K100").Value = 2590.818
$ws.Range("M100").Value = -2049.818
$ws.Range("H111").Value = 4399
$ws.Range("I111").Value = 4799.5
$ws.Range("J111").Value = 3998.5
$ws.Range("K111").Value = 14398.5
$ws.Range("L111").Value = 11995.5
$ws.Range("M111").Value = -11331.5
$ws.Range("N111").Value = -18129.5
$ws.Range("H116").Value = 4700
$ws.Range("I116").Value = 4500
$ws.Range("K116").Value = 4500
$ws.Range("M116").Value = -1058
$ws.Range("H122").Value = 1117
$ws.Range("I122").Value = 579.36365
$ws.Range("K122").Value = 1738.09095
$ws.Range("M122").Value = 711.90905
$ws.Range("H132").Value = 3496.3333
$ws.Range("I132").Value = 3496.3333
$ws.Range("K132").Value = 10488.9999
$ws.Range("M132").Value = -7958.999899999999
$ws.Range("H137").Value = 3938.4644
$ws.Range("I137").Value = 3494.88
$ws.Range("K137").Value = 10484.64
$ws.Range("M137").Value = -7934.639999999999
$ws.Range("H138").Value = 1484.3334
$ws.Range("I138").Value = 908.4286
$ws.Range("K138").Value = 2725.2858
$ws.Range("M138").Value = 2414.7142
$ws.Range("H141").Value = 10018.111
$ws.Range("I141").Value = 10018.111
$ws.Range("K141").Value = 30054.333
$ws.Range("M141").Value = -24874.333
# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 10016
$ws.Range("J21").Value = 19017
$ws.Range("L21").Value = 19017
$ws.Range("N21").Value = -19765
$ws.Range("H31").Value = 9485.5
$ws.Range("I31").Value = 9485.5
$ws.Range("K31").Value = 9485.5
$ws.Range("M31").Value = -9191.5
$ws.Range("H32").Value = 2459.0256
$ws.Range("I32").Value = 2260.6052
$ws.Range("K32").Value = 2260.6052
$ws.Range("M32").Value = -1973.6052
$ws.Range("H63").Value = 15159.583
$ws.Range("I63").Value = 14191.5
$ws.Range("K63").Value = 14191.5
$ws.Range("M63").Value = -13505.5
$ws.Range("H66").Value = 15159.583
$ws.Range("I66").Value = 14191.5
$ws.Range("K66").Value = 70957.5
$ws.Range("M66").Value = -67525.5
$ws.Range("H80").Value = 81071.766
$ws.Range("J80").Value = 81071.766
$ws.Range("L80").Value = 81071.766
$ws.Range("N80").Value = -83067.766
$ws.Range("H83").Value = 81071.766
$ws.Range("J83").Value = 81071.766
$ws.Range("L83").Value = 243215.298
$ws.Range("N83").Value = -253199.298
# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 860.4737
$ws.Range("I94").Value = 947.0909
$ws.Range("J94").Value = 741.375
$ws.Range("K94").Value = 947.0909
$ws.Range("L94").Value = 741.375
$ws.Range("M94").Value = -496.0909
$ws.Range("N94").Value = -1643.375
$ws.Range("H134").Value = 6175.9585
$ws.Range("I134").Value = 7194.375
$ws.Range("K134").Value = 21583.125
$ws.Range("M134").Value = -19048.125
# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 974.5
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H22").Value = 8000519.5
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 13333899
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 13333899
$ws.Range("M22").Value = -100
$ws.Range("N22").Value = -13334599
$ws.Range("H31").Value = 2537.3572
$ws.Range("I31").Value = 2688.6365
$ws.Range("K31").Value = 2688.6365
$ws.Range("M31").Value = -2393.6365
$ws.Range("H34").Value = 2537.3572
$ws.Range("I34").Value = 2688.6365
$ws.Range("K34").Value = 2688.6365
$ws.Range("M34").Value = -2486.6365
$ws.Range("I37").Value = 25000
$ws.Range("J37").Value = 24998
$ws.Range("K37").Value = 25000
$ws.Range("L37").Value = 24998
$ws.Range("M37").Value = -24893
$ws.Range("N37").Value = -25212
$ws.Range("H58").Value = 3490
$ws.Range("I58").Value = 2400.3333
$ws.Range("K58").Value = 2400.3333
$ws.Range("M58").Value = -2197.3333
$ws.Range("H68").Value = 34000
$ws.Range("I68").Value = 34000
$ws.Range("K68").Value = 34000
$ws.Range("M68").Value = -33251
$ws.Range("H71").Value = 34000
$ws.Range("I71").Value = 34000
$ws.Range("K71").Value = 102000
$ws.Range("M71").Value = -98256
$ws.Range("H107").Value = 376.72726
$ws.Range("I107").Value = 313.77777
$ws.Range("K107").Value = 313.77777
$ws.Range("M107").Value = 1606.22223
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1970
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 2898.0625
$ws.Range("I134").Value = 3120.7693
$ws.Range("J134").Value = 1933
$ws.Range("K134").Value = 9362.3079
$ws.Range("L134").Value = 5799
$ws.Range("M134").Value = -6827.3079
$ws.Range("N134").Value = -10869
$ws.Range("H136").Value = 3490
$ws.Range("I136").Value = 2400.3333
$ws.Range("K136").Value = 7200.999899999999
$ws.Range("M136").Value = -4650.999899999999
$ws.Range("H140").Value = 29990
$ws.Range("J140").Value = 29990
$ws.Range("L140").Value = 29990
$ws.Range("N140").Value = -40350
# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 341.14285
$ws.Range("I18").Value = 273
$ws.Range("K18").Value = 819
$ws.Range("M18").Value = -650
$ws.Range("H21").Value = 1001
$ws.Range("I21").Value = 1001
$ws.Range("K21").Value = 3003
$ws.Range("M21").Value = -2830
$ws.Range("H62").Value = 10528.5
$ws.Range("J62").Value = 10504.667
$ws.Range("L62").Value = 31514.001
$ws.Range("N62").Value = -32886.001
$ws.Range("H63").Value = 5193
$ws.Range("I63").Value = 3590.6667
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 10772.0001
$ws.Range("L63").Value = 30000
$ws.Range("M63").Value = -10023.0001
$ws.Range("N63").Value = -31498
$ws.Range("H65").Value = 10528.5
$ws.Range("J65").Value = 10504.667
$ws.Range("L65").Value = 94542.003
$ws.Range("N65").Value = -101406.003
$ws.Range("H66").Value = 5193
$ws.Range("I66").Value = 3590.6667
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 32316.0003
$ws.Range("L66").Value = 90000
$ws.Range("M66").Value = -28572.0003
$ws.Range("N66").Value = -97488
$ws.Range("H87").Value = 4997
$ws.Range("I87").Value = 4997
$ws.Range("K87").Value = 14991
$ws.Range("M87").Value = -13743
$ws.Range("H88").Value = 25000
$ws.Range("J88").Value = 25000
$ws.Range("L88").Value = 75000
$ws.Range("N88").Value = -75856
$ws.Range("H90").Value = 4997
$ws.Range("I90").Value = 4997
$ws.Range("K90").Value = 44973
$ws.Range("M90").Value = -38733
$ws.Range("H91").Value = 25000
$ws.Range("J91").Value = 25000
$ws.Range("L91").Value = 75000
$ws.Range("N91").Value = -77964
$ws.Range("H92").Value = 304.75
$ws.Range("I92").Value = 304.75
$ws.Range("K92").Value = 914.25
$ws.Range("M92").Value = 333.75
$ws.Range("H102").Value = 5800
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 5800
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 17400
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -22268
$ws.Range("H108").Value = 115.5
$ws.Range("I108").Value = 115.5
$ws.Range("K108").Value = 346.5
$ws.Range("M108").Value = 2533.5
# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1999
$ws.Range("I113").Value = 1999
$ws.Range("K113").Value = 1999
$ws.Range("M113").Value = 171
$ws.Range("H122").Value = 2801.5454
$ws.Range("I122").Value = 2805.5625
$ws.Range("J122").Value = 2790.8333
$ws.Range("K122").Value = 8416.6875
$ws.Range("L122").Value = 8372.499899999999
$ws.Range("M122").Value = -5966.6875
$ws.Range("N122").Value = -13272.4999
$ws.Range("H126").Value = 3068.5
$ws.Range("I126").Value = 2424.6667
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 7274.000100000001
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -4804.000100000001
$ws.Range("N126").Value = -19940
# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1499
$ws.Range("I61").Value = 1499
$ws.Range("K61").Value = 1499
$ws.Range("M61").Value = -1297
$ws.Range("H113").Value = 1499
$ws.Range("I113").Value = 1499
$ws.Range("K113").Value = 1499
$ws.Range("M113").Value = 671
$ws.Range("H132").Value = 2001.9231
$ws.Range("I132").Value = 1928.1875
$ws.Range("K132").Value = 5784.5625
$ws.Range("M132").Value = -3254.5625
# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9682.375
$ws.Range("I62").Value = 7000.6665
$ws.Range("K62").Value = 7000.6665
$ws.Range("M62").Value = -6376.6665
$ws.Range("H65").Value = 9682.375
$ws.Range("I65").Value = 7000.6665
$ws.Range("K65").Value = 35003.3325
$ws.Range("M65").Value = -31883.3325
$ws.Range("H96").Value = 3242
$ws.Range("J96").Value = 4181.3335
$ws.Range("L96").Value = 4181.3335
$ws.Range("N96").Value = -6927.3335
$ws.Range("H107").Value = 245
$ws.Range("I107").Value = 245
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 735
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1185
$ws.Range("N107").ClearContents()
$ws.Range("H126").Value = 3783.64
$ws.Range("I126").Value = 3334.7
$ws.Range("J126").Value = 5579.4
$ws.Range("K126").Value = 10004.1
$ws.Range("L126").Value = 16738.2
$ws.Range("M126").Value = -7534.099999999999
$ws.Range("N126").Value = -21678.2
